$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A5").Value = 521642198
$ws.Range("A6").Value = 441895023

$ws.Range("A7").Select()
